# -----------------------------------------------------------------------
# Adds 8 new "success summary" notification-template rows (for each of the
# 6 existing languages) to Sheet1, matching commit "added templates for
# summary". Columns are: A=lang_code, B=code, C=descr, D=is_active.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1605

$languages = @("eng", "ara", "fra", "hin", "kan", "tam")

$templates = @(
    @{ Code = "cust-and-down-my-card-success-summary"; Descr = "Success summary to customize and download my card" },
    @{ Code = "order-a-physical-card-success-summary"; Descr = "Success summary to order a physical card" },
    @{ Code = "share-cred-with-partner-success-summary"; Descr = "Success summary to share my credential with a partner" },
    @{ Code = "lock-unlock-auth-success-summary"; Descr = "Success summary to lock/unlock various authentication types" },
    @{ Code = "update-demo-data-success-summary"; Descr = "Success summary to self update demographic data" },
    @{ Code = "gen-or-revoke-vid-success-summary"; Descr = "Success summary to generate or revoke VIDs" },
    @{ Code = "get-my-uin-card-success-summary"; Descr = "Success summary to get my UIN card" },
    @{ Code = "verify-my-phone-email-success-summary"; Descr = "Success summary to verify my phone and email" }
)

# Pre-seed the new shared strings in the exact order they first appear in
# the authored workbook: all 8 descriptions, then all 8 codes.
$seedRow = $startRow
foreach ($t in $templates) {
    $ws.Range("C$seedRow").Value = $t.Descr
    $seedRow = $seedRow + 1
}
$seedRow = $startRow
foreach ($t in $templates) {
    $ws.Range("B$seedRow").Value = $t.Code
    $seedRow = $seedRow + 1
}

# An existing cell holding the text value "TRUE" (as opposed to a boolean),
# used as the copy source for column D so new cells keep the same shared
# string / text formatting instead of Excel auto-converting to a boolean.
$trueCell = $ws.Range("D2")

$row = $startRow
foreach ($lang in $languages) {
    foreach ($t in $templates) {
        $ws.Range("A$row").Value = $lang
        $ws.Range("B$row").Value = $t.Code
        $ws.Range("C$row").Value = $t.Descr
        $trueCell.Copy() | Out-Null
        $ws.Range("D$row").PasteSpecial(-4104) | Out-Null
        $row = $row + 1
    }
}

$excel.CutCopyMode = 0

# Restore the selection to match the authored workbook state.
$ws.Range("C1647").Select() | Out-Null
